$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$oldValues = @(
  "85-20=65",
  "93-3=90",
  "61-57=4",
  "28-17=11",
  "57+14=71",
  "6+14=20",
  "61-42=19",
  "79+18=97",
  "45+53=98",
  "30-16=14",
  "35-14=21",
  "91-50=41",
  "73-35=38",
  "96-91=5",
  "46+35=81",
  "90-85=5",
  "40-35=5",
  "10+37=47",
  "67+13=80",
  "10+34=44",
  "40-34=6",
  "1+31=32",
  "79+18=97",
  "36-24=12",
  "98-69=29",
  "45-10=35",
  "84-18=66",
  "98-87=11",
  "2+44=46",
  "57-57=0",
  "78-42=36",
  "96-46=50",
  "44+36=80",
  "64-2=62",
  "31-20=11",
  "70-28=42",
  "95-4=91",
  "47+52=99",
  "73-62=11",
  "65-39=26",
  "84-11=73",
  "3+0=3",
  "74-29=45",
  "62-25=37",
  "26+35=61",
  "41+23=64",
  "78-6=72",
  "95-77=18",
  "8+40=48",
  "56+18=74",
  "59-50=9",
  "23+52=75",
  "63+28=91",
  "27+31=58",
  "35+41=76",
  "81-5=76",
  "16+15=31",
  "49-15=34",
  "30+27=57",
  "62+7=69",
  "53+5=58",
  "8+76=84",
  "62-61=1",
  "60-44=16",
  "34+45=79",
  "67-54=13",
  "52-27=25",
  "90-16=74",
  "21+1=22",
  "45+10=55",
  "0+7=7",
  "1+52=53",
  "7+48=55",
  "39-38=1",
  "34+10=44",
  "95-94=1",
  "0+43=43",
  "62+31=93",
  "80-80=0",
  "11+77=88",
  "96-21=75",
  "38-2=36",
  "74+19=93",
  "65-21=44",
  "28+8=36",
  "35+40=75",
  "9+82=91",
  "21+16=37",
  "79+13=92",
  "75-50=25",
  "30+8=38",
  "7+27=34",
  "66+33=99",
  "64+31=95",
  "72-69=3",
  "8+44=52",
  "99-21=78",
  "95-38=57",
  "53-39=14",
  "84-9=75"
)
$values = @(
  "10+7=17",
  "21+62=83",
  "47-43=4",
  "67-56=11",
  "59+8=67",
  "89-83=6",
  "20+64=84",
  "12+70=82",
  "32+13=45",
  "35+37=72",
  "31+13=44",
  "40+1=41",
  "78-56=22",
  "22+20=42",
  "71-16=55",
  "64-29=35",
  "38-27=11",
  "29+36=65",
  "2+11=13",
  "71-34=37",
  "23+27=50",
  "6+10=16",
  "54+30=84",
  "97-32=65",
  "51-49=2",
  "60-31=29",
  "19+70=89",
  "36+28=64",
  "38+5=43",
  "69-45=24",
  "13+57=70",
  "62-46=16",
  "55-31=24",
  "41-17=24",
  "12+63=75",
  "53-38=15",
  "6+93=99",
  "0-0=0",
  "19+65=84",
  "25+64=89",
  "95-28=67",
  "57-41=16",
  "7+62=69",
  "59+23=82",
  "45+22=67",
  "43-20=23",
  "47+19=66",
  "64+1=65",
  "81-31=50",
  "72+7=79",
  "84-17=67",
  "54-9=45",
  "58-36=22",
  "94-54=40",
  "7+2=9",
  "78+11=89",
  "91-69=22",
  "93+6=99",
  "6+61=67",
  "55+25=80",
  "90-41=49",
  "47+19=66",
  "92-28=64",
  "7+2=9",
  "11+33=44",
  "79-58=21",
  "2+42=44",
  "57+42=99",
  "35+18=53",
  "92+1=93",
  "42-32=10",
  "53+42=95",
  "67-24=43",
  "65-4=61",
  "43+20=63",
  "30+54=84",
  "21+17=38",
  "87-26=61",
  "58-26=32",
  "31+22=53",
  "53+38=91",
  "95-69=26",
  "75-66=9",
  "77+20=97",
  "26+69=95",
  "79-11=68",
  "17+34=51",
  "46+0=46",
  "64-5=59",
  "35+57=92",
  "38+10=48",
  "4+28=32",
  "67+6=73",
  "40-6=34",
  "89-13=76",
  "4+44=48",
  "17+5=22",
  "23+23=46",
  "56-31=25",
  "48+12=60"
)
$idx = 0
$mismatches = 0
for ($r = 1; $r -le 20; $r++) {
  for ($c = 1; $c -le 5; $c++) {
    $cell = $t.Cell($r, $c)
    $rng = $cell.Range
    $curText = $rng.Text
    $curClean = $curText.Substring(0, $curText.Length - 2)
    $expectedOld = $oldValues[$idx]
    if ($curClean -ne $expectedOld) {
      $mismatches = $mismatches + 1
      Write-Host "Mismatch at idx" $idx ": expected [" $expectedOld "] found [" $curClean "]"
    }
    $rng.Text = $values[$idx]
    $idx = $idx + 1
  }
}
Write-Host "Updated" $idx "cells with" $mismatches "mismatches"